$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59-97 down to 60-98
$ws.Rows("59:59").Insert()

# Fill in the new row 59 with the weekly price record for Damasco
# (Castle Brite, Primera, Region de O'Higgins, $/bandeja 10 kilos)
$ws.Range("A59").Value = 9
$ws.Range("B59").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").Value = 44893
$ws.Range("D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100103
$ws.Range("H59").Value = "Frutos de hueso (carozo)"
$ws.Range("I59").Value = 100103003
$ws.Range("J59").Value = "Damasco"
$ws.Range("K59").Value = "Castle Brite"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 200
$ws.Range("N59").Value = 15000
$ws.Range("O59").Value = 15000
$ws.Range("P59").Value = 15000
$ws.Range("Q59").Value = "$/bandeja 10 kilos"
$ws.Range("R59").Value = "Región de O'Higgins"
$ws.Range("S59").Value = 1500
$ws.Range("T59").Value = 10
